{"js": "// Replace the multiplication problems in the table with the new operands,\n// matching the commit's xml diff (25 <w:t> text replacements).\nconst replacements = [\n  [\"536\u00d76=\", \"967\u00d79=\"],\n  [\"406\u00d74=\", \"426\u00d77=\"],\n  [\"555\u00d76=\", \"702\u00d76=\"],\n  [\"440\u00d79=\", \"190\u00d78=\"],\n  [\"354\u00d77=\", \"400\u00d73=\"],\n  [\"472\u00d79=\", \"614\u00d77=\"],\n  [\"251\u00d73=\", \"722\u00d72=\"],\n  [\"317\u00d72=\", \"465\u00d73=\"],\n  [\"317\u00d75=\", \"490\u00d77=\"],\n  [\"254\u00d79=\", \"403\u00d79=\"],\n  [\"315\u00d75=\", \"479\u00d74=\"],\n  [\"402\u00d77=\", \"429\u00d74=\"],\n  [\"957\u00d74=\", \"999\u00d78=\"],\n  [\"907\u00d79=\", \"388\u00d76=\"],\n  [\"684\u00d76=\", \"834\u00d76=\"],\n  [\"436\u00d74=\", \"982\u00d72=\"],\n  [\"805\u00d79=\", \"841\u00d77=\"],\n  [\"458\u00d74=\", \"916\u00d75=\"],\n  [\"997\u00d73=\", \"352\u00d73=\"],\n  [\"388\u00d75=\", \"775\u00d78=\"],\n  [\"851\u00d73=\", \"645\u00d72=\"],\n  [\"639\u00d77=\", \"933\u00d76=\"],\n  [\"109\u00d75=\", \"364\u00d73=\"],\n  [\"528\u00d72=\", \"343\u00d77=\"],\n  [\"179\u00d73=\", \"124\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the table with the new operands,\n# matching the commit's xml diff (25 text replacements via Find/Replace).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"536\u00d76=\", \"967\u00d79=\"),\n    @(\"406\u00d74=\", \"426\u00d77=\"),\n    @(\"555\u00d76=\", \"702\u00d76=\"),\n    @(\"440\u00d79=\", \"190\u00d78=\"),\n    @(\"354\u00d77=\", \"400\u00d73=\"),\n    @(\"472\u00d79=\", \"614\u00d77=\"),\n    @(\"251\u00d73=\", \"722\u00d72=\"),\n    @(\"317\u00d72=\", \"465\u00d73=\"),\n    @(\"317\u00d75=\", \"490\u00d77=\"),\n    @(\"254\u00d79=\", \"403\u00d79=\"),\n    @(\"315\u00d75=\", \"479\u00d74=\"),\n    @(\"402\u00d77=\", \"429\u00d74=\"),\n    @(\"957\u00d74=\", \"999\u00d78=\"),\n    @(\"907\u00d79=\", \"388\u00d76=\"),\n    @(\"684\u00d76=\", \"834\u00d76=\"),\n    @(\"436\u00d74=\", \"982\u00d72=\"),\n    @(\"805\u00d79=\", \"841\u00d77=\"),\n    @(\"458\u00d74=\", \"916\u00d75=\"),\n    @(\"997\u00d73=\", \"352\u00d73=\"),\n    @(\"388\u00d75=\", \"775\u00d78=\"),\n    @(\"851\u00d73=\", \"645\u00d72=\"),\n    @(\"639\u00d77=\", \"933\u00d76=\"),\n    @(\"109\u00d75=\", \"364\u00d73=\"),\n    @(\"528\u00d72=\", \"343\u00d77=\"),\n    @(\"179\u00d73=\", \"124\u00d75=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
